$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 103, shifting existing rows 103:121 down to 104:122
$ws.Rows(103).Insert()

# Populate the new row 103 with the new weekly data record.
$ws.Cells.Item(103, 1).Value = 1
$ws.Cells.Item(103, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(103, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(103, 4).Value = 45124
$ws.Cells.Item(103, 5).Value = 15
$ws.Cells.Item(103, 6).Value = 100112040
$ws.Cells.Item(103, 7).Value = "Cilantro"
$ws.Cells.Item(103, 8).Value = "Sin especificar"
$ws.Cells.Item(103, 9).Value = "Primera"
$ws.Cells.Item(103, 10).Value = 250
$ws.Cells.Item(103, 11).Value = 800
$ws.Cells.Item(103, 12).Value = 1000
$ws.Cells.Item(103, 13).Value = 900
$ws.Cells.Item(103, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(103, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(103, 16).Value = 450
$ws.Cells.Item(103, 17).Value = 2
$ws.Cells.Item(103, 18).Value = "Hortaliza"
